$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.537.37"
$ws.Range("E2").Value = "  -4.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.266.75"
$ws.Range("E3").Value = "  -6.23%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.27"
$ws.Range("E5").Value = "  -2.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.23"
$ws.Range("E6").Value = "  -4.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("E7").Value = "  -4.90%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.256.12"
$ws.Range("E9").Value = "  -6.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.603"
$ws.Range("E10").Value = "  -4.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.05"
$ws.Range("E12").Value = "  -2.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  -3.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.73"
$ws.Range("E14").Value = "  -5.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.815.84"
$ws.Range("E15").Value = "  -5.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "17.70"
$ws.Range("E16").Value = "  -4.24%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.294.19"
$ws.Range("E17").Value = "  -5.36%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.115"
$ws.Range("E18").Value = "  -4.59%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.48"
$ws.Range("E19").Value = "  -4.56%  "

$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "62.752.29"
$ws.Range("E20").Value = "  -4.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.958"
$ws.Range("E21").Value = "  -2.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "408.22"
$ws.Range("E22").Value = "  -1.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.34"
$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.97"
$ws.Range("E24").Value = "  -1.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.40"
$ws.Range("E25").Value = "  +5.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.92"
$ws.Range("E26").Value = "  -4.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.40"
$ws.Range("E27").Value = "  -3.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.68"
$ws.Range("E28").Value = "  -5.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.47"
$ws.Range("E29").Value = "  -5.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.59"
$ws.Range("E30").Value = "  -5.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.24"
$ws.Range("E31").Value = "  -3.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.18"
$ws.Range("E32").Value = "  -4.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "565.81"
$ws.Range("E33").Value = "  -7.07%  "

$ws.Range("E34").Value = "  -4.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.24"
$ws.Range("E35").Value = "  -3.73%  "

$ws.Range("E36").Value = "  -0.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.144"
$ws.Range("E37").Value = "  -1.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "34.55"
$ws.Range("E38").Value = "  -7.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.36"
$ws.Range("E39").Value = "  +3.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0724"
$ws.Range("E40").Value = "  -7.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.359"
$ws.Range("E41").Value = "  -5.23%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.087.71"
$ws.Range("E42").Value = "  -8.48%  "

$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("E44").Value = "  -0.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.72"
$ws.Range("E45").Value = "  -3.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0394"
$ws.Range("E46").Value = "  -4.37%  "

$ws.Range("E47").Value = "  -6.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.58"
$ws.Range("E48").Value = "  -4.55%  "

$ws.Range("E49").Value = "  -4.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.37"
$ws.Range("E50").Value = "  -4.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.88"
$ws.Range("E51").Value = "  -6.92%  "
